# Apply updated symbol list values (Price and Volume(1h) columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '304.91' }
    @{ Cell = 'E2'; Value = '6.21%' }
    @{ Cell = 'D3'; Value = '32.19' }
    @{ Cell = 'E3'; Value = '8.99%' }
    @{ Cell = 'D4'; Value = '5.338' }
    @{ Cell = 'E4'; Value = '4.40%' }
    @{ Cell = 'D5'; Value = '0.07329' }
    @{ Cell = 'E5'; Value = '9.32%' }
    @{ Cell = 'D6'; Value = '7.873' }
    @{ Cell = 'E6'; Value = '7.44%' }
    @{ Cell = 'D7'; Value = '3.787' }
    @{ Cell = 'E7'; Value = '11.29%' }
    @{ Cell = 'D8'; Value = '1.526' }
    @{ Cell = 'E8'; Value = '12.12%' }
    @{ Cell = 'D9'; Value = '0.9201' }
    @{ Cell = 'E9'; Value = '0.11%' }
    @{ Cell = 'D10'; Value = '0.01708' }
    @{ Cell = 'E10'; Value = '2,538.32%' }
    @{ Cell = 'D11'; Value = '0.1698' }
    @{ Cell = 'E11'; Value = '6.60%' }
    @{ Cell = 'D12'; Value = '0.07617' }
    @{ Cell = 'E12'; Value = '12.69%' }
    @{ Cell = 'D13'; Value = '0.08026' }
    @{ Cell = 'E13'; Value = '4.03%' }
    @{ Cell = 'D14'; Value = '0.03070' }
    @{ Cell = 'E14'; Value = '4.18%' }
    @{ Cell = 'D15'; Value = '0.09907' }
    @{ Cell = 'E15'; Value = '10.23%' }
    @{ Cell = 'D16'; Value = '0.001531' }
    @{ Cell = 'E16'; Value = '-4.09%' }
    @{ Cell = 'D17'; Value = '0.04564' }
    @{ Cell = 'E17'; Value = '1.59%' }
    @{ Cell = 'D18'; Value = '0.006180' }
    @{ Cell = 'E18'; Value = '-1.68%' }
    @{ Cell = 'D19'; Value = '3.462' }
    @{ Cell = 'E19'; Value = '0.25%' }
    @{ Cell = 'D20'; Value = '2.244' }
    @{ Cell = 'E20'; Value = '0.71%' }
    @{ Cell = 'D21'; Value = '0.3285' }
    @{ Cell = 'E21'; Value = '2.20%' }
    @{ Cell = 'D22'; Value = '0.1341' }
    @{ Cell = 'E22'; Value = '2.42%' }
    @{ Cell = 'D23'; Value = '4.335' }
    @{ Cell = 'E23'; Value = '6.65%' }
    @{ Cell = 'D24'; Value = '0.1634' }
    @{ Cell = 'E24'; Value = '3.34%' }
    @{ Cell = 'D25'; Value = '0.001232' }
    @{ Cell = 'E25'; Value = '3.39%' }
    @{ Cell = 'D26'; Value = '0.004485' }
    @{ Cell = 'E26'; Value = '8.84%' }
    @{ Cell = 'D27'; Value = '0.0001317' }
    @{ Cell = 'E27'; Value = '9.93%' }
    @{ Cell = 'D28'; Value = '0.0001748' }
    @{ Cell = 'E28'; Value = '8.20%' }
    @{ Cell = 'D40'; Value = '0.04552' }
    @{ Cell = 'E40'; Value = '6.61%' }
    @{ Cell = 'D41'; Value = '0.006999' }
    @{ Cell = 'E41'; Value = '3.64%' }
    @{ Cell = 'D42'; Value = '0.1353' }
    @{ Cell = 'E42'; Value = '9.14%' }
    @{ Cell = 'D43'; Value = '0.002290' }
    @{ Cell = 'E43'; Value = '2.82%' }
    @{ Cell = 'D44'; Value = '0.01432' }
    @{ Cell = 'E44'; Value = '18.81%' }
    @{ Cell = 'D45'; Value = '0.00006109' }
    @{ Cell = 'E45'; Value = '7.14%' }
    @{ Cell = 'E46'; Value = '-4.13%' }
    @{ Cell = 'D47'; Value = '0.01305' }
    @{ Cell = 'E47'; Value = '0.02%' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}

